$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 39
$ws.Range("H39").Value = 676.38464
$ws.Range("I39").Value = 268
$ws.Range("K39").Value = 804
$ws.Range("M39").Value = -508
# Row 51
$ws.Range("H51").Value = 3891.4614
$ws.Range("I51").Value = 3658.25
$ws.Range("J51").Value = 3995.111
$ws.Range("K51").Value = 3658.25
$ws.Range("L51").Value = 3995.111
$ws.Range("M51").Value = -3174.25
$ws.Range("N51").Value = -4963.111
# Row 57
$ws.Range("H57").Value = 134085.28
$ws.Range("J57").Value = 134085.28
$ws.Range("L57").Value = 402255.84
$ws.Range("N57").Value = -403253.84
# Row 64
$ws.Range("H64").Value = 4145.364
$ws.Range("I64").Value = 3599.8333
$ws.Range("J64").Value = 4800
$ws.Range("K64").Value = 3599.8333
$ws.Range("L64").Value = 4800
$ws.Range("M64").Value = -3351.8333
$ws.Range("N64").Value = -5296
# Row 67
$ws.Range("H67").Value = 4145.364
$ws.Range("I67").Value = 3599.8333
$ws.Range("J67").Value = 4800
$ws.Range("K67").Value = 3599.8333
$ws.Range("L67").Value = 4800
$ws.Range("M67").Value = -2741.8333
$ws.Range("N67").Value = -6516
# Row 70
$ws.Range("H70").Value = 1702.35
$ws.Range("I70").Value = 1941.7693
$ws.Range("J70").Value = 1257.7142
$ws.Range("K70").Value = 5825.3079
$ws.Range("L70").Value = 3773.1426
$ws.Range("M70").Value = -5555.3079
$ws.Range("N70").Value = -4313.142599999999
# Row 73
$ws.Range("H73").Value = 1702.35
$ws.Range("I73").Value = 1941.7693
$ws.Range("J73").Value = 1257.7142
$ws.Range("K73").Value = 5825.3079
$ws.Range("L73").Value = 3773.1426
$ws.Range("M73").Value = -4889.3079
$ws.Range("N73").Value = -5645.142599999999
# Row 98
$ws.Range("H98").Value = 1723.9166
$ws.Range("I98").Value = 1676.8889
$ws.Range("K98").Value = 1676.8889
$ws.Range("M98").Value = -178.8888999999999
# Row 116
$ws.Range("H116").Value = 21499.5
$ws.Range("J116").Value = 10000
$ws.Range("L116").Value = 10000
$ws.Range("N116").Value = -16884
# Row 122
$ws.Range("H122").Value = 1723.9166
$ws.Range("I122").Value = 1676.8889
$ws.Range("K122").Value = 5030.6667
$ws.Range("M122").Value = -2580.6667
# Row 129
$ws.Range("H129").Value = 3498.875
$ws.Range("I129").Value = 3284.4285
$ws.Range("K129").Value = 9853.2855
$ws.Range("M129").Value = -4853.2855
# Row 137
$ws.Range("H137").Value = 1147.1333
$ws.Range("I137").Value = 1173.3636
$ws.Range("K137").Value = 3520.0908
$ws.Range("M137").Value = -970.0907999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 116
$ws.Range("H116").Value = 12293761
$ws.Range("I116").Value = 15547631
$ws.Range("J116").Value = 1361.2222
$ws.Range("K116").Value = 15547631
$ws.Range("L116").Value = 1361.2222
$ws.Range("M116").Value = -15545337
$ws.Range("N116").Value = -5949.2222
# Row 2
$ws.Range("H2").Value = 12293761
$ws.Range("I2").Value = 15547631
$ws.Range("J2").Value = 1361.2222
$ws.Range("K2").Value = 15547631
$ws.Range("L2").Value = 1361.2222
$ws.Range("M2").Value = -15547518
$ws.Range("N2").Value = -1587.2222
# Row 45
$ws.Range("H45").Value = 19303.715
$ws.Range("J45").Value = 807
$ws.Range("L45").Value = 807
$ws.Range("N45").Value = -1561
# Row 97
$ws.Range("I97").Value = 30307344
$ws.Range("K97").Value = 30307344
$ws.Range("M97").Value = -30306848

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 9277.546
$ws.Range("I122").Value = 6277.769
$ws.Range("K122").Value = 18833.307
$ws.Range("M122").Value = -16383.307
# Row 107
$ws.Range("H107").Value = 717.5143
$ws.Range("I107").Value = 1009.3077
$ws.Range("J107").Value = 545.0909
$ws.Range("K107").Value = 3027.9231
$ws.Range("L107").Value = 1635.2727
$ws.Range("M107").Value = -1107.9231
$ws.Range("N107").Value = -5475.2727
# Row 100
$ws.Range("H100").Value = 678.26666
$ws.Range("I100").Value = 597.4545000000001
$ws.Range("J100").Value = 900.5
$ws.Range("K100").Value = 1194.909
$ws.Range("L100").Value = 1801
$ws.Range("M100").Value = -653.9090000000001
$ws.Range("N100").Value = -2883

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 137
$ws.Range("H137").Value = 8101.273
$ws.Range("I137").Value = 3511
$ws.Range("J137").Value = 13609.6
$ws.Range("K137").Value = 10533
$ws.Range("L137").Value = 40828.8
$ws.Range("M137").Value = -5433
$ws.Range("N137").Value = -51028.8
# Row 5
$ws.Range("H5").Value = 849.75
$ws.Range("I5").Value = 849.75
$ws.Range("K5").Value = 2549.25
$ws.Range("M5").Value = -2437.25
# Row 74
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
# Row 77
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
# Row 132
$ws.Range("H132").Value = 2044.2222
$ws.Range("I132").Value = 1737.375
$ws.Range("K132").Value = 15636.375
$ws.Range("M132").Value = -13106.375
# Row 135
$ws.Range("H135").Value = 849.75
$ws.Range("I135").Value = 849.75
$ws.Range("K135").Value = 7647.75
$ws.Range("M135").Value = -5112.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 388889020
$ws.Range("I2").Value = 777777800
$ws.Range("K2").Value = 777777800
$ws.Range("M2").Value = -777777687
# Row 107
$ws.Range("H107").Value = 582.4
$ws.Range("I107").Value = 539.1429000000001
$ws.Range("K107").Value = 539.1429000000001
$ws.Range("M107").Value = 1380.8571
# Row 132
$ws.Range("H132").Value = 9847.182000000001
$ws.Range("I132").Value = 10413.125
$ws.Range("K132").Value = 31239.375
$ws.Range("M132").Value = -28709.375
# Row 80
$ws.Range("H80").Value = 18405060
$ws.Range("J80").Value = 5624.1
$ws.Range("L80").Value = 5624.1
$ws.Range("N80").Value = -7620.1
# Row 83
$ws.Range("H83").Value = 18405060
$ws.Range("J83").Value = 5624.1
$ws.Range("L83").Value = 28120.5
$ws.Range("N83").Value = -38104.5
# Row 113
$ws.Range("H113").Value = 3161.25
$ws.Range("I113").Value = 2097.75
$ws.Range("K113").Value = 2097.75
$ws.Range("M113").Value = 72.25
# Row 126
$ws.Range("H126").Value = 4529.154
$ws.Range("I126").Value = 3660.375
$ws.Range("K126").Value = 10981.125
$ws.Range("M126").Value = -8511.125

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 12293761
$ws.Range("I3").Value = 15547631
$ws.Range("J3").Value = 1361.2222
$ws.Range("K3").Value = 15547631
$ws.Range("L3").Value = 1361.2222
$ws.Range("M3").Value = -15547517
$ws.Range("N3").Value = -1589.2222
# Row 22
$ws.Range("H22").Value = 228.85715
$ws.Range("I22").Value = 228.85715
$ws.Range("K22").Value = 228.85715
$ws.Range("M22").Value = -55.85714999999999
# Row 94
$ws.Range("H94").Value = 1121.8667
$ws.Range("I94").Value = 321.2
$ws.Range("J94").Value = 2723.2
$ws.Range("K94").Value = 321.2
$ws.Range("L94").Value = 2723.2
$ws.Range("M94").Value = 129.8
$ws.Range("N94").Value = -3625.2
# Row 107
$ws.Range("H107").Value = 4441.5
$ws.Range("I107").Value = 4179.5625
$ws.Range("K107").Value = 4179.5625
$ws.Range("M107").Value = -2259.5625

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 2078.5334
$ws.Range("I22").Value = 2021.25
$ws.Range("K22").Value = 2021.25
$ws.Range("M22").Value = -1671.25
# Row 107
$ws.Range("H107").Value = 1142.5264
$ws.Range("I107").Value = 643.3333
$ws.Range("J107").Value = 1998.2858
$ws.Range("K107").Value = 643.3333
$ws.Range("L107").Value = 1998.2858
$ws.Range("M107").Value = 1276.6667
$ws.Range("N107").Value = -5838.2858
# Row 58
$ws.Range("H58").Value = 6641.3184
$ws.Range("I58").Value = 7801.5386
$ws.Range("K58").Value = 7801.5386
$ws.Range("M58").Value = -7598.5386
# Row 59
$ws.Range("H59").Value = 71699.60000000001
$ws.Range("I59").Value = 49500
$ws.Range("J59").Value = 77249.5
$ws.Range("K59").Value = 49500
$ws.Range("L59").Value = 77249.5
$ws.Range("M59").Value = -48355
$ws.Range("N59").Value = -79539.5
# Row 105
$ws.Range("H105").Value = 2211
$ws.Range("I105").Value = 2211
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2211
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -464
$ws.Range("N105").ClearContents()
# Row 136
$ws.Range("H136").Value = 6641.3184
$ws.Range("I136").Value = 7801.5386
$ws.Range("K136").Value = 23404.6158
$ws.Range("M136").Value = -20854.6158

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 727.8889
$ws.Range("I22").Value = 765.3333
$ws.Range("J22").Value = 709.1667
$ws.Range("K22").Value = 765.3333
$ws.Range("L22").Value = 709.1667
$ws.Range("M22").Value = -470.3333
$ws.Range("N22").Value = -1299.1667
# Row 27
$ws.Range("H27").Value = 727.8889
$ws.Range("I27").Value = 765.3333
$ws.Range("J27").Value = 709.1667
$ws.Range("K27").Value = 765.3333
$ws.Range("L27").Value = 709.1667
$ws.Range("M27").Value = -658.3333
$ws.Range("N27").Value = -923.1667
# Row 46
$ws.Range("H46").Value = 2066.5386
$ws.Range("I46").Value = 1987
$ws.Range("J46").Value = 2116.25
$ws.Range("K46").Value = 1987
$ws.Range("L46").Value = 2116.25
$ws.Range("M46").Value = -1799
$ws.Range("N46").Value = -2492.25
